$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.617.79'
$ws.Range('E2').Value = '  +2.68%  '
# Row 3
$ws.Range('D3').Value = '2.473.90'
$ws.Range('E3').Value = '  +2.30%  '
# Row 4
$ws.Range('E4').Value = '  +0.14%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.87'
$ws.Range('E5').Value = '  +2.74%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.86'
$ws.Range('E6').Value = '  +4.29%  '
# Row 7
$ws.Range('E7').Value = '  -0.13%  '
# Row 8
$ws.Range('E8').Value = '  +2.07%  '
# Row 9
$ws.Range('E9').Value = '  +4.95%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.155'
$ws.Range('E10').Value = '  +0.80%  '
# Row 11
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.364'
$ws.Range('E11').Value = '  +4.31%  '
# Row 12
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.32'
$ws.Range('E12').Value = '  +2.54%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.29'
$ws.Range('E13').Value = '  +4.49%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000185'
$ws.Range('E14').Value = '  +7.27%  '
# Row 15
$ws.Range('D15').Value = '2.905.08'
$ws.Range('E15').Value = '  +1.92%  '
# Row 16
$ws.Range('D16').Value = '63.291.53'
$ws.Range('E16').Value = '  +2.36%  '
# Row 17
$ws.Range('D17').Value = '2.482.88'
$ws.Range('E17').Value = '  +2.81%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.58'
$ws.Range('E18').Value = '  +2.47%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.35'
$ws.Range('E19').Value = '  +7.71%  '
# Row 20
$ws.Range('E20').Value = '  +3.11%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '328.97'
$ws.Range('E21').Value = '  +1.80%  '
# Row 22
$ws.Range('B22').Value = 'SuiNetwork'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.95'
$ws.Range('E22').Value = '  +12.30%  '
# Row 23
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.23%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.51'
$ws.Range('E24').Value = '  +0.91%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '632.85'
$ws.Range('E25').Value = '  +14.41%  '
# Row 26
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.91'
$ws.Range('E26').Value = '  +2.00%  '
# Row 27
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000106'
$ws.Range('E27').Value = '  +13.83%  '
# Row 28
$ws.Range('E28').Value = '  +2.38%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.53'
$ws.Range('E29').Value = '  +10.39%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.43'
$ws.Range('E30').Value = '  +3.14%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.37%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.146'
$ws.Range('E32').Value = '  -0.40%  '
# Row 33
$ws.Range('E33').Value = '  +3.09%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.21'
$ws.Range('E34').Value = '  +10.43%  '
# Row 35
$ws.Range('E35').Value = '  +4.16%  '
# Row 36
$ws.Range('E36').Value = '  -0.14%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.387'
$ws.Range('E37').Value = '  +2.33%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.56'
$ws.Range('E38').Value = '  +2.44%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.02'
$ws.Range('E39').Value = '  +2.27%  '
# Row 40
$ws.Range('E40').Value = '  +2.82%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '146.52'
$ws.Range('E41').Value = '  -3.90%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.69'
$ws.Range('E42').Value = '  +20.99%  '
# Row 43
$ws.Range('E43').Value = '  -0.03%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '150.83'
$ws.Range('E44').Value = '  +2.66%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.78'
$ws.Range('E45').Value = '  +4.13%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0552'
$ws.Range('E46').Value = '  +4.94%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.19'
$ws.Range('E47').Value = '  +7.02%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.611'
$ws.Range('E48').Value = '  +2.75%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0241'
$ws.Range('E49').Value = '  +6.01%  '
# Row 50
$ws.Range('E50').Value = '  +1.02%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.750'
$ws.Range('E51').Value = '  +5.52%  '
